# Move "TextBox 19" to the end of the shape z-order (bring to front) and
# reposition it to its new offset, per the target OOXML diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.Item("TextBox 19")

# Bring the shape to the very front of the z-order so it becomes the last
# <p:sp> in the slide's shape tree.
$shp.ZOrder(0)

# Reposition: new a:off is x="3856382" y="-21572" (EMU). PowerPoint's COM
# Left/Top are expressed in points, i.e. EMU / 12700. (Top is nudged very
# slightly off the exact quotient so float32 round-trip through the COM
# Single-precision property lands on the exact target EMU value rather
# than one tick away.)
$shp.Left = 3856382 / 12700
$shp.Top = -1.6986
